# Apply the latest coinranking.com snapshot values scraped by the
# GitHub Actions job onto the existing "cryptos" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts plain-looking numeric text (e.g. "209.01") into a
# real number when assigned to .Value, which would round-trip it as
# 209.00999999999999 and drop the original text formatting. Prefixing
# with a literal apostrophe forces Excel to keep storing it as text,
# exactly like the source data (every Price/Volume cell is a string).
function Set-TextValue($range, [string]$value) {
    if ($value.Trim() -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

# Row 2: Bitcoin
Set-TextValue $ws.Range('D2') '25.845.21'
Set-TextValue $ws.Range('E2') '  -1.44%  '
# Row 3: Ethereum
Set-TextValue $ws.Range('D3') '1.621.98'
Set-TextValue $ws.Range('E3') '  -3.02%  '
# Row 4: TetherUSD
Set-TextValue $ws.Range('E4') '  -0.23%  '
# Row 5: BNB
Set-TextValue $ws.Range('D5') '209.01'
Set-TextValue $ws.Range('E5') '  -1.23%  '
# Row 6: XRP
Set-TextValue $ws.Range('D6') '0.5203'
Set-TextValue $ws.Range('E6') '  -1.42%  '
# Row 7: USDC
Set-TextValue $ws.Range('E7') '  -0.20%  '
# Row 8: Cardano
Set-TextValue $ws.Range('D8') '0.2564'
Set-TextValue $ws.Range('E8') '  -3.07%  '
# Row 9: Dogecoin
Set-TextValue $ws.Range('D9') '0.06264'
Set-TextValue $ws.Range('E9') '  -0.24%  '
# Row 10: Solana
Set-TextValue $ws.Range('D10') '20.40'
Set-TextValue $ws.Range('E10') '  -4.38%  '
# Row 11: TRON
Set-TextValue $ws.Range('D11') '0.07535'
Set-TextValue $ws.Range('E11') '  -0.43%  '
# Row 12: WrappedEther
Set-TextValue $ws.Range('D12') '1.638.46'
Set-TextValue $ws.Range('E12') '  -2.12%  '
# Row 13: Polkadot
Set-TextValue $ws.Range('D13') '4.375'
Set-TextValue $ws.Range('E13') '  -1.91%  '
# Row 14: WrappedliquidstakedEther2.0
Set-TextValue $ws.Range('D14') '1.843.05'
Set-TextValue $ws.Range('E14') '  -3.07%  '
# Row 15: Polygon
Set-TextValue $ws.Range('D15') '0.5459'
Set-TextValue $ws.Range('E15') '  -2.51%  '
# Row 16: ShibaInu
Set-TextValue $ws.Range('D16') '0.0₅7933'
Set-TextValue $ws.Range('E16') '  -0.92%  '
# Row 17: Litecoin
Set-TextValue $ws.Range('D17') '64.21'
Set-TextValue $ws.Range('E17') '  -4.05%  '
# Row 18: WrappedBTC
Set-TextValue $ws.Range('D18') '25.841.75'
# Row 19: Dai
Set-TextValue $ws.Range('E19') '  -0.16%  '
# Row 20: Uniswap
Set-TextValue $ws.Range('D20') '4.649'
Set-TextValue $ws.Range('E20') '  -3.40%  '
# Row 21: BitcoinCash
Set-TextValue $ws.Range('D21') '183.92'
Set-TextValue $ws.Range('E21') '  -1.96%  '
# Row 22: Avalanche
Set-TextValue $ws.Range('D22') '10.06'
Set-TextValue $ws.Range('E22') '  -3.35%  '
# Row 23: Chainlink
Set-TextValue $ws.Range('E23') '  -2.52%  '
# Row 24: BinanceUSD
Set-TextValue $ws.Range('E24') '  -0.22%  '
# Row 25: Monero
Set-TextValue $ws.Range('D25') '144.73'
Set-TextValue $ws.Range('E25') '  -3.42%  '
# Row 26: Stellar
Set-TextValue $ws.Range('D26') '0.1209'
Set-TextValue $ws.Range('E26') '  -3.78%  '
# Row 27: Cosmos
Set-TextValue $ws.Range('D27') '7.357'
Set-TextValue $ws.Range('E27') '  -2.91%  '
# Row 28: EthereumClassic
Set-TextValue $ws.Range('E28') '  -2.59%  '
# Row 29: Toncoin
Set-TextValue $ws.Range('E29') '  -0.24%  '
# Row 30: Hedera
Set-TextValue $ws.Range('E30') '  -5.73%  '
# Row 31: PancakeSwap
Set-TextValue $ws.Range('E31') '  -3.56%  '
# Row 32: InternetComputer(DFINITY)
Set-TextValue $ws.Range('D32') '3.390'
Set-TextValue $ws.Range('E32') '  -3.15%  '
# Row 33: Filecoin
Set-TextValue $ws.Range('D33') '3.349'
Set-TextValue $ws.Range('E33') '  -2.39%  '
# Row 34: LidoDAOToken
Set-TextValue $ws.Range('D34') '1.607'
Set-TextValue $ws.Range('E34') '  -1.47%  '
# Row 35: ARBITRUM
Set-TextValue $ws.Range('D35') '0.9719'
Set-TextValue $ws.Range('E35') '  -2.98%  '
# Row 36: HuobiToken
Set-TextValue $ws.Range('D36') '2.382'
Set-TextValue $ws.Range('E36') '  -1.23%  '
# Row 37: MXToken
Set-TextValue $ws.Range('D37') '2.719'
Set-TextValue $ws.Range('E37') '  -1.11%  '
# Row 38: ImmutableX
Set-TextValue $ws.Range('D38') '0.5758'
Set-TextValue $ws.Range('E38') '  -4.91%  '
# Row 39: VeChain
Set-TextValue $ws.Range('E39') '  -1.98%  '
# Row 40: TrustWalletToken
Set-TextValue $ws.Range('D40') '0.8448'
Set-TextValue $ws.Range('E40') '  -3.37%  '
# Row 41: PaxDollar
Set-TextValue $ws.Range('D41') '1.003'
Set-TextValue $ws.Range('E41') '  -0.36%  '
# Row 42: FraxShare
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D42') '1.025.41'
Set-TextValue $ws.Range('E42') '  -6.92%  '
# Row 43: Maker
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '5.658'
Set-TextValue $ws.Range('E43') '  -7.49%  '
# Row 44: Quant
Set-TextValue $ws.Range('D44') '99.39'
Set-TextValue $ws.Range('E44') '  -0.35%  '
# Row 45: RocketPoolETH
Set-TextValue $ws.Range('D45') '1.770.20'
Set-TextValue $ws.Range('E45') '  -2.92%  '
# Row 46: BabyDogeCoin
Set-TextValue $ws.Range('E46') '  -2.45%  '
# Row 47: Frax
Set-TextValue $ws.Range('D47') '0.9983'
Set-TextValue $ws.Range('E47') '  -0.81%  '
# Row 48: Aave
Set-TextValue $ws.Range('D48') '54.57'
Set-TextValue $ws.Range('E48') '  -2.66%  '
# Row 49: EnergySwap
Set-TextValue $ws.Range('D49') '7.926'
Set-TextValue $ws.Range('E49') '  -1.25%  '
# Row 50: Cronos
Set-TextValue $ws.Range('D50') '0.05153'
Set-TextValue $ws.Range('E50') '  -1.49%  '
# Row 51: Mantle
Set-TextValue $ws.Range('D51') '0.4209'
Set-TextValue $ws.Range('E51') '  -1.13%  '
